$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A1:R1 hold text values (player bio + game-log stats stored as
# text, matching the source data's inline-string representation).
# S1 holds a real number (fantasy-points share).
$ws.Range("A1:R1").NumberFormat = "@"
$ws.Range("S1").NumberFormat = "General"

$ws.Range("A1").Value = "Lazard"
$ws.Range("B1").Value = "Allen"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "2018-12-30"
$ws.Range("E1").Value = "16"
$ws.Range("F1").Value = "23.019"
$ws.Range("G1").Value = "GNB"
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = "DET"
$ws.Range("J1").Value = "L 0-31"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = "1"
$ws.Range("M1").Value = "1"
$ws.Range("N1").Value = "7"
$ws.Range("O1").Value = "7.00"
$ws.Range("P1").Value = "0"
$ws.Range("Q1").Value = "100.0%"
$ws.Range("R1").Value = "7.00"
$ws.Range("S1").Value = 0.7
